$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update progress values (column C) ---
$ws.Range("C10").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("C17").Value = 1
$ws.Range("C18").Value = 1
$ws.Range("C20").Value = 0.7
$ws.Range("C21").Value = 1

# --- Update comments (column D) ---
# Order matters: new shared-string entries are appended in the order the
# cells are written, and needs to match the order the strings were
# (re)introduced by the original author.
$ws.Range("D15").Value = "A faire avec l'ajout de nouveaux clients"
$ws.Range("D6").Value = "Review_score_evolution, RFM etc.."
$ws.Range("D7").Value = "Kmeans car le reste pas vraiment adaptée"
$ws.Range("D14").Value = "Ari_Score sur 12 mois et stabilité avec davies_bouldin "
$ws.Range("D18").Value = "Méthode du coude silhouette sample, davies_bouldin + GridSearchCV"
$ws.Range("D5").Value = "Test fait mais pas implémenté car pas nécessaire"
$ws.Range("D21").Value = "Docstring dans les fonctions et commentaire régulier"
$ws.Range("D20").Value = "Reste encore l'indentation"
$ws.Range("D17").Value = "Méthode de test cluster Kmeans et GridSearch perso pour le DBScan"
$ws.Range("D10").Value = "2 méthodes"

# --- Fix the F4 average formula (drop trailing comma) ---
$ws.Range("F4").Formula = "=AVERAGE(C5:C7,C9:C11,C13:C15,C17:C18)"

# --- Update selection / scroll position ---
[void]$ws.Range("H15").Select()
